$wb = $excel.ActiveWorkbook

# --- Rename header columns on every sheet (shared strings: County -> Province,
#     Sub-County -> District, Ward -> Subdistrict) ---
$sheetNames = @("School", "Health Care Facilities", "Train Station")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D1").Value = "Province"
    $ws.Range("E1").Value = "District"
    $ws.Range("F1").Value = "Subdistrict"
}

# --- Update per-sheet selection to C1:G1 on every sheet ---
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    [void]$ws.Range("C1:G1").Select()
}

# --- Move the active tab from "School" (sheet 1) to "Train Station" (sheet 3) ---
$wsTrain = $wb.Worksheets.Item("Train Station")
[void]$wsTrain.Activate()
[void]$wsTrain.Range("C1:G1").Select()

# --- Recolor the (auto-added) theme to match the LibreOffice palette used by
#     the fixture's theme1.xml (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink). The
#     RGB property uses the classic COM 0x00BBGGRR byte order. ---
try {
    $scheme = $wb.Theme.ThemeColorScheme
    $scheme.Colors(1).RGB  = 0x000000   # dk1      -> 000000
    $scheme.Colors(2).RGB  = 0xFFFFFF   # lt1      -> ffffff
    $scheme.Colors(3).RGB  = 0x000000   # dk2      -> 000000
    $scheme.Colors(4).RGB  = 0xFFFFFF   # lt2      -> ffffff
    $scheme.Colors(5).RGB  = 0x03A318   # accent1  -> 18a303
    $scheme.Colors(6).RGB  = 0xA36903   # accent2  -> 0369a3
    $scheme.Colors(7).RGB  = 0x033EA3   # accent3  -> a33e03
    $scheme.Colors(8).RGB  = 0xA3038E   # accent4  -> 8e03a3
    $scheme.Colors(9).RGB  = 0x009CC9   # accent5  -> c99c00
    $scheme.Colors(10).RGB = 0x1E21C9   # accent6  -> c9211e
    $scheme.Colors(11).RGB = 0xEE0000   # hlink    -> 0000ee
    $scheme.Colors(12).RGB = 0x8B1A55   # folHlink -> 551a8b
} catch {
    Write-Output "theme recolor skipped: $_"
}
